# BFI, MI, 241111 modified 5
# Append 10 new survey response rows (438-447) to the Form_Responses1 table,
# matching the source workbook's alternating row-style banding and
# refreshing the frozen-pane selection to the new bottom of the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$tbl = $ws.ListObjects.Item(1)

# New response rows, in column order A..Y (timestamp, email, dept, student id,
# name, then the 20 survey answer columns).
$rows = @(
    @(45610.957850115738, 'jerryterryharry@gmail.com', '빅데이터', 20205162, '문진영', '1. 0-2일', '5. 6-7일', '3. 3-5일', '3. 3-5일', '3. 3-5일', '1. 0-2일', '1. 0-2일', '3. 3-5일', '5. 6-7일', '1. 0-2일', '3. 가끔', '1. 예', '1. 예', '1. 예', '5.아니오', '5.아니오', '1. 예', '3. 가끔', '3. 가끔', '5.아니오'),
    @(45610.959684016205, 'tngusvhs@gmail.com', '생명과학과', 20243529, '이수현', '1. 0-2일', '3. 3-5일', '5. 6-7일', '1. 0-2일', '5. 6-7일', '5. 6-7일', '5. 6-7일', '5. 6-7일', '3. 3-5일', '1. 0-2일', '1. 예', '3. 가끔', '1. 예', '1. 예', '1. 예', '1. 예', '1. 예', '5.아니오', '3. 가끔', '1. 예'),
    @(45610.973043310187, 'jb9517asd@naver.com', '소프트웨어학부', 20245109, '곽우주', '5. 6-7일', '5. 6-7일', '5. 6-7일', '3. 3-5일', '3. 3-5일', '1. 0-2일', '1. 0-2일', '5. 6-7일', '5. 6-7일', '5. 6-7일', '5.아니오', '3. 가끔', '3. 가끔', '5.아니오', '5.아니오', '5.아니오', '5.아니오', '5.아니오', '5.아니오', '5.아니오'),
    @(45611.03734920139, 'lhw2565@gmail.com', '미디어스쿨', 20242565, '이혜원', '1. 0-2일', '5. 6-7일', '3. 3-5일', '1. 0-2일', '3. 3-5일', '3. 3-5일', '1. 0-2일', '3. 3-5일', '3. 3-5일', '1. 0-2일', '1. 예', '1. 예', '1. 예', '3. 가끔', '3. 가끔', '1. 예', '1. 예', '3. 가끔', '1. 예', '1. 예'),
    @(45611.388949548607, 'bigeyejimmy1@naver.com', '경영학과', 20182850, '김현준', '1. 0-2일', '3. 3-5일', '1. 0-2일', '1. 0-2일', '1. 0-2일', '1. 0-2일', '1. 0-2일', '3. 3-5일', '1. 0-2일', '1. 0-2일', '3. 가끔', '3. 가끔', '3. 가끔', '5.아니오', '5.아니오', '3. 가끔', '3. 가끔', '3. 가끔', '3. 가끔', '3. 가끔'),
    @(45611.463201238425, 'yhh323@naver.com', '체육', 20184132, '유형호', '1. 0-2일', '1. 0-2일', '1. 0-2일', '1. 0-2일', '1. 0-2일', '1. 0-2일', '1. 0-2일', '3. 3-5일', '1. 0-2일', '1. 0-2일', '1. 예', '1. 예', '3. 가끔', '3. 가끔', '3. 가끔', '3. 가끔', '3. 가끔', '3. 가끔', '3. 가끔', '3. 가끔'),
    @(45611.653322129627, 'hyj13223@naver.com', '정치행정학과', 20212432, '이현진', '5. 6-7일', '5. 6-7일', '5. 6-7일', '3. 3-5일', '5. 6-7일', '1. 0-2일', '1. 0-2일', '3. 3-5일', '1. 0-2일', '5. 6-7일', '3. 가끔', '3. 가끔', '5.아니오', '5.아니오', '5.아니오', '3. 가끔', '3. 가끔', '1. 예', '5.아니오', '5.아니오'),
    @(45611.692850821761, 'chaecjb@naver.com', '디지털미디어콘텐츠전공', 20203046, '채희수', '5. 6-7일', '5. 6-7일', '5. 6-7일', '5. 6-7일', '5. 6-7일', '3. 3-5일', '1. 0-2일', '5. 6-7일', '5. 6-7일', '1. 0-2일', '1. 예', '3. 가끔', '3. 가끔', '3. 가끔', '3. 가끔', '3. 가끔', '1. 예', '3. 가끔', '5.아니오', '5.아니오'),
    @(45611.698145231479, 'bcy1976@naver.com', '빅데이터학과', 20235180, '변치윤', '1. 0-2일', '5. 6-7일', '3. 3-5일', '3. 3-5일', '3. 3-5일', '3. 3-5일', '3. 3-5일', '3. 3-5일', '3. 3-5일', '1. 0-2일', '3. 가끔', '1. 예', '1. 예', '3. 가끔', '3. 가끔', '3. 가끔', '3. 가끔', '1. 예', '5.아니오', '1. 예'),
    @(45611.705875486106, 'emf1811@naver.com', '바이오메디컬학과', 20233605, '김들', '3. 3-5일', '5. 6-7일', '3. 3-5일', '3. 3-5일', '3. 3-5일', '5. 6-7일', '1. 0-2일', '5. 6-7일', '3. 3-5일', '5. 6-7일', '3. 가끔', '1. 예', '1. 예', '1. 예', '1. 예', '5.아니오', '5.아니오', '5.아니오', '1. 예', '3. 가끔')
)

# Rows already in the sheet use a 2-row alternating style band (odd slot /
# even slot), with the very last row of the table carrying a distinct
# "final row" style. Grab one exemplar of each so new rows reuse the same
# style ids instead of generating fresh ones.
$lastRowIndex = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
$evenBandSrc = $ws.Range("A" + ($lastRowIndex - 1) + ":Y" + ($lastRowIndex - 1))
$oddBandSrc  = $ws.Range("A" + ($lastRowIndex - 2) + ":Y" + ($lastRowIndex - 2))
$finalBandSrc = $ws.Range("A" + $lastRowIndex + ":Y" + $lastRowIndex)

$xlPasteFormats = -4122

$firstNewRow = $lastRowIndex + 1
$count = $rows.Count

for ($i = 0; $i -lt $count; $i++) {
    $newRow = $tbl.ListRows.Add()
    $r = $firstNewRow + $i
    $isLast = ($i -eq ($count - 1))

    $destRange = $ws.Range("A" + $r + ":Y" + $r)
    if ($isLast) {
        $finalBandSrc.Copy()
    } elseif ($i % 2 -eq 0) {
        $evenBandSrc.Copy()
    } else {
        $oddBandSrc.Copy()
    }
    $destRange.PasteSpecial($xlPasteFormats)
    $excel.CutCopyMode = 0

    $ws.Rows.Item($r).RowHeight = 15.75

    $rowData = $rows[$i]
    $ws.Cells.Item($r, 1).Value = $rowData[0]
    $ws.Cells.Item($r, 2).Value = $rowData[1]
    $ws.Cells.Item($r, 3).Value = $rowData[2]
    $ws.Cells.Item($r, 4).Value = $rowData[3]
    $ws.Cells.Item($r, 5).Value = $rowData[4]
    for ($c = 6; $c -le 25; $c++) {
        $ws.Cells.Item($r, $c).Value = $rowData[$c - 1]
    }
}

$lastNewRow = $firstNewRow + $count - 1

# Refresh the frozen-pane view to match the new bottom-of-data selection.
$ws.Range("E" + ($lastNewRow + 6)).Select()
